$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B25").Value = "2019.5.27  18：30-21：30"
$ws.Range("C25").Value = "service层的实体管理器和方法的实现与测试"
$ws.Range("B26").Value = "2019.5.28  13：44-"

$ws.Range("B26").Select()
